$d = $word.ActiveDocument

# wdYellow = 7, wdRed = 6
$wdYellow = 7
$wdRed = 6

# Headings that simply need their text (run + paragraph mark) highlighted.
$yellowHeadings = @(
    "Promise All",
    "Promise Race",
    "Async Function with Await",
    "Async Function with Error Handling",
    "Chained Promises with Async/Await",
    "Quiz",
    "Simple Stopwatch",
    "Simple Text Adventure Game"
)

$redHeadings = @(
    "Countdown Timer"
)

function Get-HeadingParagraph($doc, $headingText) {
    $paras = $doc.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        $t = $p.Range.Text.TrimEnd([char]13)
        if ($t -eq $headingText) {
            return $p
        }
    }
    return $null
}

foreach ($h in $yellowHeadings) {
    $p = Get-HeadingParagraph $d $h
    if ($p -ne $null) {
        $p.Range.Font.HighlightColorIndex = $wdYellow
    }
}

foreach ($h in $redHeadings) {
    $p = Get-HeadingParagraph $d $h
    if ($p -ne $null) {
        $p.Range.Font.HighlightColorIndex = $wdRed
    }
}

# "Promise with Multiple Handlers" also gets highlighted and is followed by
# a new, empty (Normal-style) paragraph whose paragraph mark carries the
# same yellow highlight.
$target = Get-HeadingParagraph $d "Promise with Multiple Handlers"
if ($target -ne $null) {
    $target.Range.Font.HighlightColorIndex = $wdYellow

    $insertAt = $d.Range($target.Range.End, $target.Range.End)
    $newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $insertAt.InsertXML($newParaXml)
}
